$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 8400
$ws.Range("I18").Value = 8400
$ws.Range("K18").Value = 8400
$ws.Range("M18").Value = -8116

$ws.Range("H135").Value = 5488.077
$ws.Range("I135").Value = 1571
$ws.Range("K135").Value = 14139
$ws.Range("M135").Value = -11604

$ws.Range("H137").Value = 8430.950000000001
$ws.Range("I137").Value = 3906.5881
$ws.Range("K137").Value = 11719.7643
$ws.Range("M137").Value = -9169.764299999999

$ws.Range("H138").Value = 5362.8706
$ws.Range("J138").Value = 5385.61
$ws.Range("L138").Value = 16156.83
$ws.Range("N138").Value = -26436.83

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 133.125
$ws.Range("I4").Value = 133.125
$ws.Range("K4").Value = 133.125
$ws.Range("M4").Value = -17.125

$ws.Range("H32").Value = 8485186
$ws.Range("I32").Value = 9266918
$ws.Range("J32").Value = 42482.6
$ws.Range("K32").Value = 9266918
$ws.Range("L32").Value = 42482.6
$ws.Range("M32").Value = -9266631
$ws.Range("N32").Value = -43056.6

$ws.Range("H74").Value = 31266250
$ws.Range("I74").Value = 125001000
$ws.Range("K74").Value = 125001000
$ws.Range("M74").Value = -125000126

$ws.Range("H77").Value = 31266250
$ws.Range("I77").Value = 125001000
$ws.Range("K77").Value = 625005000
$ws.Range("M77").Value = -625000632

$ws.Range("H132").Value = 4276.173
$ws.Range("I132").Value = 2044.075
$ws.Range("K132").Value = 6132.225
$ws.Range("M132").Value = -3602.225

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1448
$ws.Range("I86").Value = 1216.8667
$ws.Range("K86").Value = 1216.8667
$ws.Range("M86").Value = -93.86670000000004

$ws.Range("H89").Value = 1448
$ws.Range("I89").Value = 1216.8667
$ws.Range("K89").Value = 6084.333500000001
$ws.Range("M89").Value = -468.3335000000006

$ws.Range("H94").Value = 1629.9412
$ws.Range("I94").Value = 1475.0834
$ws.Range("K94").Value = 1475.0834
$ws.Range("M94").Value = -1024.0834

$ws.Range("H131").Value = 15000
$ws.Range("J131").Value = 15000
$ws.Range("L131").Value = 15000
$ws.Range("N131").Value = -25080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 40000

$ws.Range("H51").Value = 21750
$ws.Range("I51").Value = 21750
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 21750
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -21014
$ws.Range("N51").ClearContents()

$ws.Range("H60").Value = 74999
$ws.Range("I60").Value = 74999
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 74999
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -74488
$ws.Range("N60").ClearContents()

$ws.Range("H61").Value = 21750
$ws.Range("I61").Value = 21750
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 21750
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -21402
$ws.Range("N61").ClearContents()

$ws.Range("H62").Value = 3999
$ws.Range("J62").Value = 3999
$ws.Range("L62").Value = 3999
$ws.Range("N62").Value = -5247

$ws.Range("H65").Value = 3999
$ws.Range("J65").Value = 3999
$ws.Range("L65").Value = 19995
$ws.Range("N65").Value = -26235

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1249.8
$ws.Range("J34").Value = 2624.75
$ws.Range("L34").Value = 7874.25
$ws.Range("N34").Value = -8042.25

$ws.Range("H39").Value = 89583.164
$ws.Range("I39").Value = 8600
$ws.Range("J39").Value = 147428.28
$ws.Range("K39").Value = 25800
$ws.Range("L39").Value = 442284.84
$ws.Range("M39").Value = -25506
$ws.Range("N39").Value = -442872.84

$ws.Range("H55").Value = 6482.3335
$ws.Range("I55").Value = 7223.75
$ws.Range("K55").Value = 21671.25
$ws.Range("M55").Value = -21494.25

$ws.Range("H68").Value = 2087.2068
$ws.Range("I68").Value = 2261.5
$ws.Range("J68").Value = 2041.7391
$ws.Range("K68").Value = 6784.5
$ws.Range("L68").Value = 6125.2173
$ws.Range("M68").Value = -5973.5
$ws.Range("N68").Value = -7747.2173

$ws.Range("H71").Value = 2087.2068
$ws.Range("I71").Value = 2261.5
$ws.Range("J71").Value = 2041.7391
$ws.Range("K71").Value = 20353.5
$ws.Range("L71").Value = 18375.6519
$ws.Range("M71").Value = -16297.5
$ws.Range("N71").Value = -26487.6519

$ws.Range("H107").Value = 743.61536
$ws.Range("I107").Value = 655.5454999999999
$ws.Range("J107").Value = 808.2
$ws.Range("K107").Value = 1966.6365
$ws.Range("L107").Value = 2424.6
$ws.Range("M107").Value = -46.63649999999984
$ws.Range("N107").Value = -6264.6

$ws.Range("H113").Value = 967.86206
$ws.Range("I113").Value = 754.2857
$ws.Range("J113").Value = 1035.8182
$ws.Range("K113").Value = 2262.8571
$ws.Range("L113").Value = 3107.4546
$ws.Range("M113").Value = -92.85710000000017
$ws.Range("N113").Value = -7447.4546

$ws.Range("H131").Value = 3563.2856
$ws.Range("I131").Value = 3767.8948
$ws.Range("K131").Value = 11303.6844
$ws.Range("M131").Value = -6263.6844

$ws.Range("H140").Value = 177294
$ws.Range("I140").Value = 200880.14
$ws.Range("K140").Value = 602640.42
$ws.Range("M140").Value = -597460.42

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8501.166999999999
$ws.Range("I70").Value = 7999.3335
$ws.Range("K70").Value = 7999.3335
$ws.Range("M70").Value = -7729.3335

$ws.Range("H73").Value = 8501.166999999999
$ws.Range("I73").Value = 7999.3335
$ws.Range("K73").Value = 7999.3335
$ws.Range("M73").Value = -7063.3335

$ws.Range("H99").Value = 28317.334
$ws.Range("J99").Value = 54812.668
$ws.Range("L99").Value = 54812.668
$ws.Range("N99").Value = -59304.668

$ws.Range("H100").Value = 117000
$ws.Range("J100").Value = 117000
$ws.Range("L100").Value = 117000
$ws.Range("N100").Value = -119164

$ws.Range("H102").Value = 2423.8064
$ws.Range("I102").Value = 1973.52
$ws.Range("K102").Value = 1973.52
$ws.Range("M102").Value = -351.52

$ws.Range("H126").Value = 5440.25
$ws.Range("I126").Value = 6056.8
$ws.Range("K126").Value = 18170.4
$ws.Range("M126").Value = -15700.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3395.36
$ws.Range("I16").Value = 3522.158
$ws.Range("J16").Value = 2993.8333
$ws.Range("K16").Value = 3522.158
$ws.Range("L16").Value = 2993.8333
$ws.Range("M16").Value = -3352.158
$ws.Range("N16").Value = -3333.8333

$ws.Range("H40").Value = 2968.739
$ws.Range("I40").Value = 2541.1052
$ws.Range("K40").Value = 2541.1052
$ws.Range("M40").Value = -2405.1052

$ws.Range("H68").Value = 907.6
$ws.Range("J68").Value = 437
$ws.Range("L68").Value = 437
$ws.Range("N68").Value = -1935

$ws.Range("H71").Value = 907.6
$ws.Range("J71").Value = 437
$ws.Range("L71").Value = 2185
$ws.Range("N71").Value = -9673

$ws.Range("H136").Value = 35764.477
$ws.Range("I136").Value = 4821.875
$ws.Range("J136").Value = 82178.375
$ws.Range("K136").Value = 14465.625
$ws.Range("L136").Value = 246535.125
$ws.Range("M136").Value = -11915.625
$ws.Range("N136").Value = -251635.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 50000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 50000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 50000
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -54992

$ws.Range("H102").Value = 134000
$ws.Range("J102").Value = 134000
$ws.Range("L102").Value = 134000
$ws.Range("N102").Value = -140490

$ws.Range("H113").Value = 1325.5
$ws.Range("I113").Value = 1340.6
$ws.Range("J113").Value = 1250
$ws.Range("K113").Value = 4021.8
$ws.Range("L113").Value = 3750
$ws.Range("M113").Value = -1851.8
$ws.Range("N113").Value = -8090

$ws.Range("H124").Value = 85214.5
$ws.Range("J124").Value = 85214.5
$ws.Range("L124").Value = 85214.5
$ws.Range("N124").Value = -95034.5

$ws.Range("H136").Value = 3336.4285
$ws.Range("I136").Value = 2269.7
$ws.Range("K136").Value = 6809.099999999999
$ws.Range("M136").Value = -4259.099999999999
